$wb = $excel.ActiveWorkbook

# This script applies updated market-price snapshot values (columns H-N)
# for the rows affected by the latest scheduled data refresh, across the
# eight job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 910.7778
$ws.Range("I9").Value = 819.8
$ws.Range("J9").Value = 1024.5
$ws.Range("K9").Value = 819.8
$ws.Range("L9").Value = 1024.5
$ws.Range("M9").Value = -650.8
$ws.Range("N9").Value = -1362.5
$ws.Range("H41").Value = 2108.2
$ws.Range("I41").Value = 265
$ws.Range("J41").Value = 2569
$ws.Range("K41").Value = 265
$ws.Range("L41").Value = 2569
$ws.Range("M41").Value = 175
$ws.Range("N41").Value = -3449
$ws.Range("H70").Value = 4108.9688
$ws.Range("I70").Value = 2086.9167
$ws.Range("J70").Value = 5322.2
$ws.Range("K70").Value = 6260.750100000001
$ws.Range("L70").Value = 15966.6
$ws.Range("M70").Value = -5990.750100000001
$ws.Range("N70").Value = -16506.6
$ws.Range("H73").Value = 4108.9688
$ws.Range("I73").Value = 2086.9167
$ws.Range("J73").Value = 5322.2
$ws.Range("K73").Value = 6260.750100000001
$ws.Range("L73").Value = 15966.6
$ws.Range("M73").Value = -5324.750100000001
$ws.Range("N73").Value = -17838.6
$ws.Range("H137").Value = 253115.16
$ws.Range("I137").Value = 266121.25
$ws.Range("K137").Value = 798363.75
$ws.Range("M137").Value = -795813.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 540
$ws.Range("I5").Value = 486.81818
$ws.Range("K5").Value = 486.81818
$ws.Range("M5").Value = -374.81818
$ws.Range("H32").Value = 16736.12
$ws.Range("I32").Value = 17740.318
$ws.Range("K32").Value = 17740.318
$ws.Range("M32").Value = -17453.318
$ws.Range("H61").Value = 3291.2856
$ws.Range("I61").Value = 2339.3
$ws.Range("K61").Value = 2339.3
$ws.Range("M61").Value = -2127.3
$ws.Range("H74").Value = 11907334
$ws.Range("I74").Value = 25001340
$ws.Range("K74").Value = 25001340
$ws.Range("M74").Value = -25000466
$ws.Range("H77").Value = 11907334
$ws.Range("I77").Value = 25001340
$ws.Range("K77").Value = 125006700
$ws.Range("M77").Value = -125002332
$ws.Range("H132").Value = 13239.288
$ws.Range("I132").Value = 15726.846
$ws.Range("J132").Value = 5776.615
$ws.Range("K132").Value = 47180.538
$ws.Range("L132").Value = 17329.845
$ws.Range("M132").Value = -44650.538
$ws.Range("N132").Value = -22389.845
$ws.Range("H136").Value = 3291.2856
$ws.Range("I136").Value = 2339.3
$ws.Range("K136").Value = 7017.900000000001
$ws.Range("M136").Value = -4467.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 540
$ws.Range("I4").Value = 486.81818
$ws.Range("K4").Value = 486.81818
$ws.Range("M4").Value = -371.81818
$ws.Range("H20").Value = 2317.9092
$ws.Range("I20").Value = 2362.5
$ws.Range("K20").Value = 2362.5
$ws.Range("M20").Value = -2115.5
$ws.Range("H22").Value = 732.3333
$ws.Range("I22").Value = 99.5
$ws.Range("J22").Value = 1998
$ws.Range("K22").Value = 99.5
$ws.Range("L22").Value = 1998
$ws.Range("M22").Value = 73.5
$ws.Range("N22").Value = -2344
$ws.Range("H86").Value = 2308.1667
$ws.Range("I86").Value = 2308.1667
$ws.Range("K86").Value = 2308.1667
$ws.Range("M86").Value = -1185.1667
$ws.Range("H89").Value = 2308.1667
$ws.Range("I89").Value = 2308.1667
$ws.Range("K89").Value = 11540.8335
$ws.Range("M89").Value = -5924.833500000001
$ws.Range("H99").Value = 1995.8334
$ws.Range("I99").Value = 1795
$ws.Range("K99").Value = 1795
$ws.Range("M99").Value = -297
$ws.Range("H107").Value = 2180.879
$ws.Range("I107").Value = 2100.2144
$ws.Range("K107").Value = 2100.2144
$ws.Range("M107").Value = -180.2143999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 604.2963
$ws.Range("I22").Value = 674.6087
$ws.Range("K22").Value = 674.6087
$ws.Range("M22").Value = -324.6087
$ws.Range("H31").Value = 25003848
$ws.Range("I31").Value = 35717556
$ws.Range("K31").Value = 35717556
$ws.Range("M31").Value = -35717261
$ws.Range("H34").Value = 25003848
$ws.Range("I34").Value = 35717556
$ws.Range("K34").Value = 35717556
$ws.Range("M34").Value = -35717354
$ws.Range("H58").Value = 502048.8
$ws.Range("I58").Value = 1828
$ws.Range("K58").Value = 1828
$ws.Range("M58").Value = -1625
$ws.Range("H62").Value = 33760
$ws.Range("I62").Value = 3600
$ws.Range("J62").Value = 79000
$ws.Range("K62").Value = 3600
$ws.Range("L62").Value = 79000
$ws.Range("M62").Value = -2976
$ws.Range("N62").Value = -80248
$ws.Range("H65").Value = 33760
$ws.Range("I65").Value = 3600
$ws.Range("J65").Value = 79000
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 395000
$ws.Range("M65").Value = -14880
$ws.Range("N65").Value = -401240
$ws.Range("H132").Value = 37048170
$ws.Range("I132").Value = 47629070
$ws.Range("J132").Value = 15002
$ws.Range("K132").Value = 142887210
$ws.Range("L132").Value = 45006
$ws.Range("M132").Value = -142884680
$ws.Range("N132").Value = -50066
$ws.Range("H134").Value = 1607.0256
$ws.Range("I134").Value = 1588
$ws.Range("K134").Value = 4764
$ws.Range("M134").Value = -2229
$ws.Range("H136").Value = 502048.8
$ws.Range("I136").Value = 1828
$ws.Range("K136").Value = 5484
$ws.Range("M136").Value = -2934

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 178.21428
$ws.Range("I40").Value = 90.59999999999999
$ws.Range("J40").Value = 226.88889
$ws.Range("K40").Value = 362.4
$ws.Range("L40").Value = 907.55556
$ws.Range("M40").Value = -293.4
$ws.Range("N40").Value = -1045.55556
$ws.Range("H56").Value = 52639492
$ws.Range("I56").Value = 52639492
$ws.Range("K56").Value = 52639492
$ws.Range("M56").Value = -52638962
$ws.Range("H107").Value = 601.25
$ws.Range("I107").Value = 343.5
$ws.Range("K107").Value = 1030.5
$ws.Range("M107").Value = 889.5
$ws.Range("H108").Value = 13749.5
$ws.Range("I108").Value = 1999
$ws.Range("J108").Value = 17666.334
$ws.Range("K108").Value = 5997
$ws.Range("L108").Value = 52999.00199999999
$ws.Range("M108").Value = -3117
$ws.Range("N108").Value = -58759.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 18325
$ws.Range("I44").Value = 14987.5
$ws.Range("K44").Value = 14987.5
$ws.Range("M44").Value = -14391.5
$ws.Range("H70").Value = 7483.5454
$ws.Range("I70").Value = 7396.5
$ws.Range("K70").Value = 7396.5
$ws.Range("M70").Value = -7126.5
$ws.Range("H73").Value = 7483.5454
$ws.Range("I73").Value = 7396.5
$ws.Range("K73").Value = 7396.5
$ws.Range("M73").Value = -6460.5
$ws.Range("H80").Value = 128266.445
$ws.Range("I80").Value = 147771.14
$ws.Range("K80").Value = 147771.14
$ws.Range("M80").Value = -146773.14
$ws.Range("H83").Value = 128266.445
$ws.Range("I83").Value = 147771.14
$ws.Range("K83").Value = 738855.7000000001
$ws.Range("M83").Value = -733863.7000000001
$ws.Range("H93").Value = 29500
$ws.Range("J93").Value = 29500
$ws.Range("L93").Value = 29500
$ws.Range("N93").Value = -33244
$ws.Range("H132").Value = 295900.12
$ws.Range("I132").Value = 61684.883
$ws.Range("J132").Value = 1433517
$ws.Range("K132").Value = 185054.649
$ws.Range("L132").Value = 4300551
$ws.Range("M132").Value = -182524.649
$ws.Range("N132").Value = -4305611

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4034
$ws.Range("I82").Value = 1801
$ws.Range("K82").Value = 1801
$ws.Range("M82").Value = -1440
$ws.Range("H85").Value = 4034
$ws.Range("I85").Value = 1801
$ws.Range("K85").Value = 1801
$ws.Range("M85").Value = -553
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H132").Value = 2095.81
$ws.Range("I132").Value = 2104.0752
$ws.Range("J132").Value = 1986
$ws.Range("K132").Value = 6312.225600000001
$ws.Range("L132").Value = 5958
$ws.Range("M132").Value = -3782.225600000001
$ws.Range("N132").Value = -11018
$ws.Range("H136").Value = 3739.7632
$ws.Range("I136").Value = 2538.1853
$ws.Range("K136").Value = 7614.5559
$ws.Range("M136").Value = -5064.5559
$ws.Range("H139").Value = 79399
$ws.Range("J139").Value = 79399
$ws.Range("L139").Value = 79399
$ws.Range("N139").Value = -89679

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4384.857
$ws.Range("I81").Value = 1000.25
$ws.Range("J81").Value = 8897.666999999999
$ws.Range("K81").Value = 2000.5
$ws.Range("L81").Value = 17795.334
$ws.Range("M81").Value = -939.5
$ws.Range("N81").Value = -19917.334
$ws.Range("H84").Value = 4384.857
$ws.Range("I84").Value = 1000.25
$ws.Range("J84").Value = 8897.666999999999
$ws.Range("K84").Value = 10002.5
$ws.Range("L84").Value = 88976.67
$ws.Range("M84").Value = -4698.5
$ws.Range("N84").Value = -99584.67
$ws.Range("H109").Value = 59992
$ws.Range("J109").Value = 59992
$ws.Range("L109").Value = 59992
$ws.Range("N109").Value = -62766
$ws.Range("H132").Value = 2384.6875
$ws.Range("J132").Value = 10288.25
$ws.Range("L132").Value = 30864.75
$ws.Range("N132").Value = -35924.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N118").ClearContents()
